$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.497298240661621
$ws.Range("B1").Value = 3.986066579818726
$ws.Range("C1").Value = 3.615141868591309
$ws.Range("D1").Value = 1.513768434524536
$ws.Range("E1").Value = 0.9817589521408081
